# "add high level bill" - extend the Kill-Bill tables with levels 21-30.
#
# Adds 10 new rows (levels 21-30) to the "data", "cn" and "en" sheets,
# refreshes the view selections that Excel recorded after the edit, and
# restores the original row-order for the header merge cells on "data".

$wb = $excel.ActiveWorkbook

# Best-effort cosmetic nudge of the saved window position (harmless if the
# host doesn't round-trip it).
try { $excel.ActiveWindow.Left = 320 } catch {}

$wsData = $wb.Worksheets.Item("data")
$wsCn   = $wb.Worksheets.Item("cn")
$wsEn   = $wb.Worksheets.Item("en")

# ---------------------------------------------------------------------
# "data" sheet: rows 23-32 (ids 10021-10030, levels 21-30)
# ---------------------------------------------------------------------
$dataRows = @(
    @{ Row=23; A=10021; B=39; C=5; D="Bill's Army 21"; E=15000; H=180000; I=3; J=4; K=1; L=$null; M=1; N=40000 },
    @{ Row=24; A=10022; B=42; C=5; D="Bill's Army 22"; E=15500; H=200000; I=3; J=9; K=1; L=$null; M=1; N=42000 },
    @{ Row=25; A=10023; B=45; C=5; D="Bill's Army 23"; E=16000; H=220000; I=1; J=3; K=1; L=3;     M=1; N=44000 },
    @{ Row=26; A=10024; B=48; C=5; D="Bill's Army 24"; E=16500; H=240000; I=3; J=5; K=1; L=$null; M=1; N=46000 },
    @{ Row=27; A=10025; B=51; C=5; D="Bill's Army 25"; E=17000; H=260000; I=1; J=4; K=1; L=3;     M=1; N=48000 },
    @{ Row=28; A=10026; B=54; C=5; D="Bill's Army 26"; E=17500; H=280000; I=3; J=6; K=1; L=$null; M=1; N=50000 },
    @{ Row=29; A=10027; B=57; C=5; D="Bill's Army 27"; E=18000; H=300000; I=3; J=7; K=1; L=$null; M=1; N=52000 },
    @{ Row=30; A=10028; B=60; C=5; D="Bill's Army 28"; E=18500; H=320000; I=1; J=5; K=1; L=3;     M=1; N=54000 },
    @{ Row=31; A=10029; B=63; C=5; D="Bill's Army 29"; E=19000; H=340000; I=3; J=8; K=1; L=$null; M=1; N=56000 },
    @{ Row=32; A=10030; B=66; C=5; D="Bill's Army 30"; E=19500; H=360000; I=1; J=6; K=1; L=3;     M=1; N=58000 }
)

foreach ($r in $dataRows) {
    $row = $r.Row

    $wsData.Cells.Item($row, 1).Value = $r.A          # A
    $wsData.Cells.Item($row, 2).Value = $r.B          # B
    $wsData.Cells.Item($row, 3).Value = $r.C          # C
    $wsData.Cells.Item($row, 4).Value = $r.D          # D (name, string)
    $wsData.Cells.Item($row, 5).Value = $r.E          # E (gold)
    $wsData.Cells.Item($row, 8).Value = $r.H          # H
    $wsData.Cells.Item($row, 9).Value = $r.I          # I
    $wsData.Cells.Item($row, 10).Value = $r.J         # J
    $wsData.Cells.Item($row, 11).Value = $r.K         # K
    if ($r.L -ne $null) {
        $wsData.Cells.Item($row, 12).Value = $r.L     # L (only when present)
    }
    $wsData.Cells.Item($row, 13).Value = $r.M         # M
    $wsData.Cells.Item($row, 14).Value = $r.N         # N

    # Formatting: every populated cell is centred; column E additionally
    # carries the thousands-separator number format. Only touch L when it
    # actually holds a value so empty cells stay absent from the row (as
    # in the authored edit) rather than materialising as blank styled
    # cells.
    $rng = $wsData.Range("A" + $row + ":D" + $row)
    $rng.HorizontalAlignment = -4108
    $eRng = $wsData.Range("E" + $row)
    $eRng.NumberFormat = "#,##0"
    $eRng.HorizontalAlignment = -4108
    $hk = $wsData.Range("H" + $row + ":K" + $row)
    $hk.HorizontalAlignment = -4108
    if ($r.L -ne $null) {
        $wsData.Range("L" + $row).HorizontalAlignment = -4108
    }
    $mn = $wsData.Range("M" + $row + ":N" + $row)
    $mn.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------
# "cn" sheet: rows 22-31
# ---------------------------------------------------------------------
$cnNames = @(
    "击败Bill 21","击败Bill 22","击败Bill 23","击败Bill 24","击败Bill 25",
    "击败Bill 26","击败Bill 27","击败Bill 28","击败Bill 29","击败Bill 30"
)
for ($i = 0; $i -lt 10; $i++) {
    $row = 22 + $i
    $wsCn.Cells.Item($row, 1).Value = 10021 + $i
    $wsCn.Cells.Item($row, 2).Value = $cnNames[$i]
    $cnRng = $wsCn.Range("A" + $row + ":B" + $row)
    $cnRng.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------
# "en" sheet: rows 22-31
# ---------------------------------------------------------------------
$enNames = @(
    "Bill's Army 21","Bill's Army 22","Bill's Army 23","Bill's Army 24","Bill's Army 25",
    "Bill's Army 26","Bill's Army 27","Bill's Army 28","Bill's Army 29","Bill's Army 30"
)
for ($i = 0; $i -lt 10; $i++) {
    $row = 22 + $i
    $wsEn.Cells.Item($row, 1).Value = 10021 + $i
    $wsEn.Cells.Item($row, 2).Value = $enNames[$i]
    $enRng = $wsEn.Range("A" + $row + ":B" + $row)
    $enRng.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------
# Restore the merged-header cell order on "data" (Excel re-emits these in
# a different order once the sheet grows).
# ---------------------------------------------------------------------
$wsData.Range("A1:O2").UnMerge()
$wsData.Range("I1:K1").Merge()
$wsData.Range("L1:L2").Merge()
$wsData.Range("M1:M2").Merge()
$wsData.Range("N1:N2").Merge()
$wsData.Range("O1:O2").Merge()
$wsData.Range("A1:A2").Merge()
$wsData.Range("H1:H2").Merge()
$wsData.Range("E1:G1").Merge()
$wsData.Range("D1:D2").Merge()
$wsData.Range("C1:C2").Merge()
$wsData.Range("B1:B2").Merge()

# ---------------------------------------------------------------------
# View selections recorded by Excel after the edit.
# ---------------------------------------------------------------------
$wsCn.Activate()
$wsCn.Range("B18:B31").Select()

$wsEn.Activate()
$wsEn.Range("E32").Select()

$wsData.Activate()
$wsData.Range("P29").Select()
